$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows($r1, $r2) {
    # Columns F:V hold the match-specific data (teams, scores, odds, timestamps, url).
    # Columns A:E (index, country, tournament, season, match date) stay put.
    $range1 = $ws.Range("F$r1" + ":V$r1")
    $range2 = $ws.Range("F$r2" + ":V$r2")
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

# Re-order pairs of matches that were scraped/sorted differently.
Swap-MatchRows 23 24
Swap-MatchRows 74 75
Swap-MatchRows 77 78
Swap-MatchRows 88 89

# Append two new match rows (99 and 100), copying formatting from the last
# existing data row (98) so the Indice (A) and data_partida (E) columns keep
# the same styling/number format as the rest of the sheet.
$ws.Range("A98:V98").Copy()
$ws.Range("A99:V100").PasteSpecial(-4122)

$ws.Range("A99").Value = 98
$ws.Range("B99").Value = "portugal"
$ws.Range("C99").Value = "liga-portugal-2"
$ws.Range("D99").Value = "2023-2024"
$ws.Range("E99").Value = 45255.5
$ws.Range("F99").Value = "FC Porto B"
$ws.Range("G99").Value = 3
$ws.Range("H99").Value = "Academico Viseu"
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 2.09
$ws.Range("K99").Value = "17/11/2023 12:42"
$ws.Range("L99").Value = 2.35
$ws.Range("M99").Value = "25/11/2023 11:53"
$ws.Range("N99").Value = 3.55
$ws.Range("O99").Value = "17/11/2023 12:42"
$ws.Range("P99").Value = 3.54
$ws.Range("Q99").Value = "25/11/2023 11:51"
$ws.Range("R99").Value = 3.53
$ws.Range("S99").Value = "17/11/2023 12:42"
$ws.Range("T99").Value = 3.08
$ws.Range("U99").Value = "25/11/2023 11:53"
$ws.Range("V99").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-academico-viseu/dUQKXqLj/"

$ws.Range("A100").Value = 99
$ws.Range("B100").Value = "portugal"
$ws.Range("C100").Value = "liga-portugal-2"
$ws.Range("D100").Value = "2023-2024"
$ws.Range("E100").Value = 45255.79166666666
$ws.Range("F100").Value = "Benfica B"
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = "Oliveirense"
$ws.Range("I100").Value = 1
$ws.Range("J100").Value = 1.93
$ws.Range("K100").Value = "19/11/2023 14:48"
$ws.Range("L100").Value = 2.03
$ws.Range("M100").Value = "25/11/2023 18:56"
$ws.Range("N100").Value = 3.74
$ws.Range("O100").Value = "19/11/2023 14:48"
$ws.Range("P100").Value = 3.66
$ws.Range("Q100").Value = "25/11/2023 18:56"
$ws.Range("R100").Value = 3.86
$ws.Range("S100").Value = "19/11/2023 14:48"
$ws.Range("T100").Value = 3.72
$ws.Range("U100").Value = "25/11/2023 18:56"
$ws.Range("V100").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/benfica-oliveirense/pj3Bf7S3/"
